$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.540661215782166
$ws.Range("B1").Value = 2.403077602386475
$ws.Range("C1").Value = 5.389420032501221
$ws.Range("D1").Value = 3.577815771102905
$ws.Range("E1").Value = 0.9780349135398865
